# Insert a new data row at row 686 (shifts existing rows 686..716 down to 687..717)
# and populate it with the new price-report record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(686).Insert()

$ws.Cells.Item(686, 1).Value = 5
$ws.Cells.Item(686, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(686, 3).Value = "Maule"
$ws.Cells.Item(686, 4).Value = 45147
$ws.Cells.Item(686, 5).Value = 7
$ws.Cells.Item(686, 6).Value = 100112043
$ws.Cells.Item(686, 7).Value = "Pepino ensalada"
$ws.Cells.Item(686, 8).Value = "Sin especificar"
$ws.Cells.Item(686, 9).Value = "Primera"
$ws.Cells.Item(686, 10).Value = 300
$ws.Cells.Item(686, 11).Value = 10000
$ws.Cells.Item(686, 12).Value = 10000
$ws.Cells.Item(686, 13).Value = 10000
$ws.Cells.Item(686, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(686, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(686, 16).Value = 167
$ws.Cells.Item(686, 17).Value = 60
$ws.Cells.Item(686, 18).Value = "Hortaliza"
